$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-01-06 Monday" "2025-01-07 Tuesday"

Replace-Text "33×99=" "97×47="
Replace-Text "88×67=" "37×45="
Replace-Text "75×62=" "43×49="
Replace-Text "70×44=" "34×48="
Replace-Text "64×96=" "18×57="

Replace-Text "75×53=" "75×49="
Replace-Text "81×80=" "20×19="
Replace-Text "61×23=" "51×70="
Replace-Text "38×75=" "86×34="
Replace-Text "13×82=" "54×52="

Replace-Text "55×85=" "56×24="
Replace-Text "57×15=" "59×87="
Replace-Text "35×49=" "73×56="
Replace-Text "67×91=" "61×37="
Replace-Text "43×81=" "52×41="

Replace-Text "46×17=" "69×26="
Replace-Text "66×60=" "47×70="
Replace-Text "25×57=" "33×34="
Replace-Text "99×26=" "78×19="
Replace-Text "38×84=" "69×69="

Replace-Text "19×32=" "23×92="
Replace-Text "78×17=" "86×56="
Replace-Text "31×34=" "14×96="
Replace-Text "30×79=" "27×51="
Replace-Text "85×56=" "40×79="
